{"js": "// The document contains three short paragraphs whose text reads\n// \"<id>p143r_1</id>\", \"<id>p143r_2</id>\", \"<id>p143r_3</id>\" \u2014 each one\n// originally split across three runs: \"<id>\" (Courier New formatting),\n// the bare id value (plain formatting), and \"</id>\" (Courier New\n// formatting again). The edit collapses each trio into a single run\n// that carries the opening run's formatting and the full concatenated\n// text, e.g. \"<id>p143r_1</id>\".\n//\n// Re-inserting the identical, already-concatenated text via\n// Range.insertText(..., \"Replace\") over a range that spans exactly\n// those three runs reproduces this merge: Word collapses the replaced\n// span into one run using the formatting of the first run it touches,\n// which is exactly the desired outcome.\nconst ids = [\"p143r_1\", \"p143r_2\", \"p143r_3\"];\nconst body = context.document.body;\n\nfor (const id of ids) {\n  const fullText = \"<id>\" + id + \"</id>\";\n  // matchCase keeps this exact; there is also a \"fig_p143r_1\" elsewhere\n  // in the document but searching for the full \"<id>...</id>\" wrapper\n  // text makes sure we never touch it.\n  const results = body.search(fullText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(fullText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# The document has three short paragraphs whose rendered text is\n# \"<id>p143r_1</id>\", \"<id>p143r_2</id>\" and \"<id>p143r_3</id>\".\n# Each one is originally split across three runs:\n#   run1 = \"<id>\"            (Courier New, color 7f6000, sz 18 - \"tag\" look)\n#   run2 = the bare id value (color 000000 - \"plain\" look)\n#   run3 = \"</id>\"           (same formatting as run1)\n# followed by a trailing empty run.\n#\n# The edit merges run1+run2+run3 into a single run that keeps run1's\n# formatting/identity and whose text is the full \"<id>...</id>\" string\n# (the trailing empty run is left untouched).\n#\n# To reproduce that precisely (including keeping run1's own XML\n# attributes such as w:rsidDel/w:rsidR/w:rsidRPr and the\n# xml:space=\"preserve\" on <w:t>), we delete the characters belonging to\n# run2+run3 and then append their combined text back onto run1's range\n# via InsertAfter - Word folds the appended text into run1 instead of\n# minting a brand-new run.\n\n$d = $word.ActiveDocument\n$ids = @(\"p143r_1\", \"p143r_2\", \"p143r_3\")\n$prefix = \"<id>\"\n$suffix = \"</id>\"\n\nforeach ($id in $ids) {\n    $target = $prefix + $id + $suffix\n\n    $paragraphCount = $d.Paragraphs.Count\n    for ($i = 1; $i -le $paragraphCount; $i++) {\n        $p = $d.Paragraphs($i)\n        # Paragraph.Range.Text includes the trailing paragraph mark (CR).\n        $paraText = $p.Range.Text.TrimEnd([char]13)\n\n        if ($paraText -eq $target) {\n            $paraRange = $p.Range\n            $start = $paraRange.Start\n            $prefixLen = $prefix.Length\n            $targetEnd = $start + $target.Length\n\n            # Remove everything after \"<id>\" (i.e. the id value + \"</id>\").\n            $tailRange = $d.Range($start + $prefixLen, $targetEnd)\n            $tailRange.Delete()\n\n            # Re-append the removed text onto the \"<id>\" run so it gets\n            # absorbed into that same run instead of creating a new one.\n            $headRange = $d.Range($start, $start + $prefixLen)\n            $headRange.InsertAfter($target.Substring($prefixLen))\n\n            break\n        }\n    }\n}\n"}
